$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Successful run for 0.5s timestep - updated results row
$ws.Range("A2").Value = 281.5
$ws.Range("B2").Value = 236.5
$ws.Range("C2").Value = 16.315770574572472
$ws.Range("D2").Value = 33.395457983595726
$ws.Range("E2").Value = 15.727193057959482
$ws.Range("F2").Value = 30.588891713187067

# Column widths refreshed to fit the new (wider) numbers
$ws.Columns.Item(1).ColumnWidth = 20
$ws.Columns.Item(2).ColumnWidth = 18.666666666666668
$ws.Columns.Item(3).ColumnWidth = 19.166666666666668
$ws.Columns.Item(4).ColumnWidth = 22.166666666666668
$ws.Columns.Item(5).ColumnWidth = 17.833333333333336
$ws.Columns.Item(6).ColumnWidth = 20.833333333333336
